$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.855.77'
$ws.Range('E2').Value = '  +0.57%  '

$ws.Range('D3').Value = '1.641.72'
$ws.Range('E3').Value = '  +0.52%  '

$ws.Range('E4').Value = '  -0.67%  '

$ws.Range('D5').Value = '''216.54'
$ws.Range('E5').Value = '  -0.53%  '

$ws.Range('D6').Value = '''0.507'
$ws.Range('E6').Value = '  +1.97%  '

$ws.Range('E7').Value = '  -0.63%  '

$ws.Range('E8').Value = '  +1.91%  '

$ws.Range('D9').Value = '''0.0621'
$ws.Range('E9').Value = '  +0.32%  '

$ws.Range('D10').Value = '''19.90'
$ws.Range('E10').Value = '  +4.91%  '

$ws.Range('E11').Value = '  +0.31%  '

$ws.Range('D12').Value = '1.870.76'
$ws.Range('E12').Value = '  +0.50%  '

$ws.Range('D13').Value = '1.640.85'
$ws.Range('E13').Value = '  +0.59%  '

$ws.Range('E14').Value = '  +0.70%  '

$ws.Range('D15').Value = '''0.528'
$ws.Range('E15').Value = '  +1.41%  '

$ws.Range('D16').Value = '''66.43'
$ws.Range('E16').Value = '  +3.88%  '

$ws.Range('D17').Value = '26.859.98'
$ws.Range('E17').Value = '  +0.68%  '

$ws.Range('D18').Value = '0.0₃0730'
$ws.Range('E18').Value = '  +1.51%  '

$ws.Range('D19').Value = '''218.71'
$ws.Range('E19').Value = '  +3.70%  '

$ws.Range('E20').Value = '  -0.60%  '

$ws.Range('D21').Value = '''6.67'
$ws.Range('E21').Value = '  +8.37%  '

$ws.Range('E22').Value = '  +1.84%  '

$ws.Range('D23').Value = '''2.41'
$ws.Range('E23').Value = '  +4.41%  '

$ws.Range('E24').Value = '  +0.13%  '

$ws.Range('D25').Value = '''145.93'
$ws.Range('E25').Value = '  -0.53%  '

$ws.Range('E26').Value = '  -0.76%  '

$ws.Range('D27').Value = '''7.38'
$ws.Range('E27').Value = '  +5.35%  '

$ws.Range('E28').Value = '  +1.60%  '

$ws.Range('E29').Value = '  +2.07%  '

$ws.Range('E30').Value = '  +1.88%  '

$ws.Range('E31').Value = '  -0.22%  '

$ws.Range('E32').Value = '  +0.45%  '

$ws.Range('D33').Value = '''2.98'
$ws.Range('E33').Value = '  +1.89%  '

$ws.Range('E34').Value = '  +3.03%  '

$ws.Range('D35').Value = '''2.45'
$ws.Range('E35').Value = '  +0.16%  '

$ws.Range('D36').Value = '1.236.72'
$ws.Range('E36').Value = '  -1.81%  '

$ws.Range('E37').Value = '  +1.36%  '

$ws.Range('E38').Value = '  +3.70%  '

$ws.Range('D39').Value = '''0.834'
$ws.Range('E39').Value = '  +4.98%  '

$ws.Range('E40').Value = '  -0.59%  '

$ws.Range('E41').Value = '  +0.96%  '

$ws.Range('D42').Value = '''5.37'
$ws.Range('E42').Value = '  +2.58%  '

$ws.Range('D43').Value = '1.782.37'
$ws.Range('E43').Value = '  +0.61%  '

$ws.Range('D44').Value = '''2.09'
$ws.Range('E44').Value = '  -3.33%  '

$ws.Range('D45').Value = '''60.88'
$ws.Range('E45').Value = '  +2.18%  '

$ws.Range('D46').Value = '''91.52'
$ws.Range('E46').Value = '  +0.59%  '

$ws.Range('E47').Value = '  +1.10%  '

$ws.Range('B48').Value = 'Cronos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D48').Value = '''0.0514'
$ws.Range('E48').Value = '  -0.52%  '

$ws.Range('B49').Value = 'Algorand'
$ws.Range('C49').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D49').Value = '''0.0972'
$ws.Range('E49').Value = '  +2.15%  '

$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').Value = '''7.57'
$ws.Range('E50').Value = '  +2.03%  '

$ws.Range('B51').Value = 'Mantle'
$ws.Range('C51').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D51').Value = '''0.405'
$ws.Range('E51').Value = '  -0.07%  '
